# "Generate Report for Handback" - refresh the localization-status report:
#   * Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     for every locale (Overview + each per-locale sheet).
#   * Per-locale "Latest Handback DateTime" is refreshed to the new handback time.
#   * The stale "handback file is not latest" Error Detail is cleared now that
#     the handback is in sync.
#   * A couple of report columns are widened/narrowed to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.0846354166667
$overview.Columns.Item(6).ColumnWidth = 29.0846354166667

# --- zh-cn sheet ------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-21 02:53:12"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.0846354166667
$zhcn.Columns.Item(16).ColumnWidth = 12.7513020833333

# --- de-de sheet ------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-21 02:53:19"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.0846354166667
$dede.Columns.Item(16).ColumnWidth = 12.7513020833333

Write-Host "Report regenerated for handback."
